$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.661.97'
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.917.92'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.64'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.37%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4927'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2971'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06751'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.884.60'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.15'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07358'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.178'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.64'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6699'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.651.32'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007928'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.50'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.38%  '
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.138.60'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.362'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +11.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '203.75'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +7.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.312'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.648'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.24'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.75'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.954'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.483'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +5.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.376'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09169'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.063'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05267'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7412'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.116'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.51%  '
$ws.Range("E36").Value = '  -1.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01845'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.719'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9258'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.36%  '
$ws.Range("E40").Value = '  -2.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4459'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.85'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +26.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.986'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '106.28'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.002'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1390'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.631'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.31'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +5.02%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.999'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05880'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4028'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.22%  '
